$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates: Team_ID (col A) and Project_Name (col B) for rows 2-6 ---
$ws.Range("A2").Value = 277
$ws.Range("B2").Value = "P128 - USYMGDL new site initialisation"

$ws.Range("A3").Value = 278
$ws.Range("B3").Value = "P1005 - P TUR TUN"

$ws.Range("A4").Value = 279
$ws.Range("B4").Value = "P1007 - P ESP UAE"

$ws.Range("A5").Value = 280
$ws.Range("B5").Value = "P1005 - P TUR TUN"

$ws.Range("A6").Value = 281
$ws.Range("B6").Value = "P1007 - P SAF UAE"

# --- New font for the newly created work-hierarchy entry (row 2) ---
$ws.Range("B2").Font.Name = "OpenSans"
$ws.Range("B2").Font.Color = 0
$ws.Range("B2").Font.Size = 11

# --- Old trailing rows (7-9) no longer hold data, but keep their formatting ---
$ws.Range("A7:B9").ClearContents()

# --- Normalise the thin box-border look (drop the stray "apply fill" flag some
#     of these cells carried) so borders render consistently across the block ---
foreach ($addr in "A5","B5","A6","B7","B8") {
  $cell = $ws.Range($addr)
  $cell.Borders.Color = 0
  $cell.Borders.LineStyle = 1
}

# --- Column B is now much wider to fit the longer project names ---
$ws.Columns("B").ColumnWidth = 42.95

# --- Rows 16 and 17 are fully cleared (content + formatting) and drop out of the sheet ---
$ws.Range("A16:B17").Clear()

# --- Final active selection left on B9 ---
$ws.Range("B9").Select()
